$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 276.02
$ws.Range("I15").Value = 276.02
$ws.Range("K15").Value = 828.0599999999999
$ws.Range("M15").Value = -659.0599999999999
$ws.Range("H32").Value = 715.6667
$ws.Range("I32").Value = 563.6667
$ws.Range("J32").Value = 867.6667
$ws.Range("K32").Value = 563.6667
$ws.Range("L32").Value = 867.6667
$ws.Range("M32").Value = -237.6667
$ws.Range("N32").Value = -1519.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1139.4865
$ws.Range("I2").Value = 1003.8571
$ws.Range("K2").Value = 1003.8571
$ws.Range("M2").Value = -890.8570999999999
$ws.Range("H32").Value = 4832.33
$ws.Range("I32").Value = 4075.9895
$ws.Range("J32").Value = 19202.8
$ws.Range("K32").Value = 4075.9895
$ws.Range("L32").Value = 19202.8
$ws.Range("M32").Value = -3788.9895
$ws.Range("N32").Value = -19776.8
$ws.Range("H61").Value = 261648.47
$ws.Range("I61").Value = 5899.148
$ws.Range("K61").Value = 5899.148
$ws.Range("M61").Value = -5687.148
$ws.Range("H97").Value = 1575.5
$ws.Range("I97").Value = 1655.091
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 1655.091
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = -1159.091
$ws.Range("N97").Value = -1692
$ws.Range("H116").Value = 1139.4865
$ws.Range("I116").Value = 1003.8571
$ws.Range("K116").Value = 1003.8571
$ws.Range("M116").Value = 1290.1429
$ws.Range("H122").Value = 1167344.1
$ws.Range("I122").Value = 1167344.1
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3502032.3
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3499582.3
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 261648.47
$ws.Range("I136").Value = 5899.148
$ws.Range("K136").Value = 17697.444
$ws.Range("M136").Value = -15147.444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1139.4865
$ws.Range("I3").Value = 1003.8571
$ws.Range("K3").Value = 1003.8571
$ws.Range("M3").Value = -889.8570999999999
$ws.Range("H94").Value = 1859.7646
$ws.Range("I94").Value = 1183.6
$ws.Range("J94").Value = 2825.7144
$ws.Range("K94").Value = 1183.6
$ws.Range("L94").Value = 2825.7144
$ws.Range("M94").Value = -732.5999999999999
$ws.Range("N94").Value = -3727.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 223856.17
$ws.Range("I31").Value = 1696.0613
$ws.Range("J31").Value = 742229.75
$ws.Range("K31").Value = 1696.0613
$ws.Range("L31").Value = 742229.75
$ws.Range("M31").Value = -1401.0613
$ws.Range("N31").Value = -742819.75
$ws.Range("H34").Value = 223856.17
$ws.Range("I34").Value = 1696.0613
$ws.Range("J34").Value = 742229.75
$ws.Range("K34").Value = 1696.0613
$ws.Range("L34").Value = 742229.75
$ws.Range("M34").Value = -1494.0613
$ws.Range("N34").Value = -742633.75
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1864.5555
$ws.Range("I92").Value = 795
$ws.Range("J92").Value = 2720.2
$ws.Range("K92").Value = 2385
$ws.Range("L92").Value = 8160.599999999999
$ws.Range("M92").Value = -1137
$ws.Range("N92").Value = -10656.6
$ws.Range("H131").Value = 2942175.2
$ws.Range("I131").Value = 10000670
$ws.Range("J131").Value = 1135.8334
$ws.Range("K131").Value = 30002010
$ws.Range("L131").Value = 3407.5002
$ws.Range("M131").Value = -29996970
$ws.Range("N131").Value = -13487.5002
$ws.Range("H132").Value = 982022.9
$ws.Range("I132").Value = 640.8
$ws.Range("J132").Value = 1309150.2
$ws.Range("K132").Value = 5767.2
$ws.Range("L132").Value = 11782351.8
$ws.Range("M132").Value = -3237.2
$ws.Range("N132").Value = -11787411.8
$ws.Range("H140").Value = 7847.143
$ws.Range("I140").Value = 5373.846
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 16121.538
$ws.Range("L140").Value = 120000
$ws.Range("M140").Value = -10941.538
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 20000
$ws.Range("K58").Value = 20000
$ws.Range("M58").Value = -19723
$ws.Range("H111").Value = 28646.25
$ws.Range("J111").Value = 28646.25
$ws.Range("L111").Value = 28646.25
$ws.Range("N111").Value = -34780.25
$ws.Range("H122").Value = 46298784
$ws.Range("I122").Value = 59158556
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 177475668
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -177473218
$ws.Range("N122").Value = -15700
$ws.Range("H123").Value = 13245
$ws.Range("J123").Value = 13245
$ws.Range("L123").Value = 13245
$ws.Range("N123").Value = -18145
$ws.Range("H132").Value = 4844.625
$ws.Range("I132").Value = 5734.2905
$ws.Range("J132").Value = 3222.2942
$ws.Range("K132").Value = 17202.8715
$ws.Range("L132").Value = 9666.882599999999
$ws.Range("M132").Value = -14672.8715
$ws.Range("N132").Value = -14726.8826

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 9980
$ws.Range("J57").Value = 9980
$ws.Range("L57").Value = 9980
$ws.Range("N57").Value = -11112
$ws.Range("H93").Value = 1074.875
$ws.Range("I93").Value = 1014.1429
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1014.1429
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 233.8570999999999
$ws.Range("N93").Value = -3996
$ws.Range("H122").Value = 3408670
$ws.Range("I122").Value = 4209533.5
$ws.Range("K122").Value = 12628600.5
$ws.Range("M122").Value = -12626150.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3830.2856
$ws.Range("I122").Value = 2499.5
$ws.Range("J122").Value = 4362.6
$ws.Range("K122").Value = 7498.5
$ws.Range("L122").Value = 13087.8
$ws.Range("M122").Value = -5048.5
$ws.Range("N122").Value = -17987.8
$ws.Range("H136").Value = 2431.7454
$ws.Range("I136").Value = 2470.516
$ws.Range("K136").Value = 7411.548000000001
$ws.Range("M136").Value = -4861.548000000001
